$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add the new "algorithm" paragraphs right after the "Algorithm"
#    title paragraph (the last paragraph in the body).
# ---------------------------------------------------------------------

# Make a fresh blank paragraph at the very end of the body to use as an
# insertion anchor, then drop the new content in just before it via
# InsertXML (raw WordprocessingML paragraphs, so none of them inherit
# the "Title" style/pPr of the paragraph they're being added after).
$lastRange = $d.Paragraphs.Last.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

$anchorPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertion = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)

$newBodyXml =
  '<w:p/>' +
  '<w:p><w:r><w:t>Ask user for variable (miles traveled)</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Ask user for variable (MPG)</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Ask user for variable (gas cost)</w:t></w:r></w:p>' +
  '<w:p/>' +
  '<w:p><w:r><w:t>Calculat</w:t></w:r><w:r><w:t>e total gas cost of trip</w:t></w:r></w:p>' +
  '<w:p/>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  $newBodyXml +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertion.InsertXML($packageXml)

# InsertXML merges its very last fragment paragraph into whatever
# paragraph used to sit at the insertion point (the blank anchor we just
# made, plus the trailing protector "<w:p/>" from $newBodyXml both land
# after our real content) - drop the two leftover placeholder paragraphs
# so "Calculate total gas cost of trip" ends up directly before sectPr.
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()
$tail = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($tail.Range.Start - 1, $tail.Range.End).Delete()

# ---------------------------------------------------------------------
# 2) Namespace / mc:Ignorable touch-ups that a newer Word build stamps
#    onto every part it resaves (w16du everywhere; oel only where the
#    "o" namespace is already declared), plus the numbering durableId.
# ---------------------------------------------------------------------
$full = $d.Content
$pkg = $full.WordOpenXML

$pkg = $pkg.Replace(
  'xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"',
  'xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'
)

$pkg = $pkg.Replace(
  'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash"',
  'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16du="http://schemas.microsoft.com/office/word/2023/wordml/word16du" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash"'
)

$pkg = $pkg.Replace(
  'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"',
  'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh w16du wp14"'
)

$pkg = $pkg.Replace(
  'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh"',
  'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh w16du"'
)

$pkg = $pkg.Replace(
  '<w:num w:numId="1">',
  '<w:num w:numId="1" w16cid:durableId="335308480">'
)

$full.WordOpenXML = $pkg
